$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "description" header in column C
$ws.Range("C1").Value = "description"

# Descriptions for the existing parameter rows
$ws.Range("C2").Value = "choose the population variant from input file"
$ws.Range("C3").Value = "choose scenario(s), comma spereated. Or choose all"

# New parameter row: restauration_building_type bias
$ws.Range("A4").Value = "restauration_building_type bias"
$ws.Range("B4").Value = "no"
$ws.Range("C4").Value = "choose how building types are affected by restauration: Either from the number of buildings (no)  or from the bias given in input file "

# Column widths (closest achievable values given the runtime's pixel-based width quantization)
$ws.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws.Columns.Item(2).ColumnWidth = 19.833333333333332
$ws.Columns.Item(3).ColumnWidth = 54.166666666666664

# Match the active selection recorded in the saved workbook
$ws.Range("C10").Select()
